$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            $t = $tr.Text
            if ($t -eq "os::rtos::Thread") {
                $tr.Characters(1, $t.Length).Text = "os::rtos::thread"
            }
            elseif ($t -eq "(Old API)") {
                $tr.Characters(1, $t.Length).Text = "(ARM API)"
            }
        }
    }
}
